$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.891087532043457
$ws.Range("B1").Value = 1.682713866233826
$ws.Range("C1").Value = 4.176331996917725
$ws.Range("D1").Value = 3.600562334060669
$ws.Range("E1").Value = 0.4175925254821777
